$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.030041153708001
$ws.Cells.Item(2, 4).Value = 1.040561006008125
$ws.Cells.Item(2, 5).Value = 1.029795012411994
$ws.Cells.Item(2, 6).Value = 1.051298934084347
$ws.Cells.Item(2, 9).Value = 1.037348532801631
$ws.Cells.Item(2, 10).Value = 1.03518459254792
$ws.Cells.Item(2, 11).Value = 1.043343005981962
$ws.Cells.Item(2, 12).Value = 1.03260788800209
$ws.Cells.Item(2, 13).Value = 1.054050822061252
$ws.Cells.Item(2, 14).Value = 1.036654673406671
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.031038494719884
$ws.Cells.Item(3, 4).Value = 1.041353055835113
$ws.Cells.Item(3, 5).Value = 1.030643014995211
$ws.Cells.Item(3, 6).Value = 1.052246401872528
$ws.Cells.Item(3, 9).Value = 1.037574045068964
$ws.Cells.Item(3, 10).Value = 1.035822975720883
$ws.Cells.Item(3, 11).Value = 1.043945552700262
$ws.Cells.Item(3, 12).Value = 1.033264003739008
$ws.Cells.Item(3, 13).Value = 1.054810561386808
$ws.Cells.Item(3, 14).Value = 1.037293963156962
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.031684092038213
$ws.Cells.Item(4, 4).Value = 1.041865413140153
$ws.Cells.Item(4, 5).Value = 1.031192330605814
$ws.Cells.Item(4, 6).Value = 1.052859602128239
$ws.Cells.Item(4, 9).Value = 1.037718256148537
$ws.Cells.Item(4, 10).Value = 1.036235721652386
$ws.Cells.Item(4, 11).Value = 1.044334634746721
$ws.Cells.Item(4, 12).Value = 1.033688514516682
$ws.Cells.Item(4, 13).Value = 1.055301665217269
$ws.Cells.Item(4, 14).Value = 1.037707295235031
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.031955560741703
$ws.Cells.Item(5, 4).Value = 1.042080770689293
$ws.Cells.Item(5, 5).Value = 1.031423405815694
$ws.Cells.Item(5, 6).Value = 1.053117420422063
$ws.Cells.Item(5, 9).Value = 1.037778472446761
$ws.Cells.Item(5, 10).Value = 1.036409160176023
$ws.Cells.Item(5, 11).Value = 1.044498011037773
$ws.Cells.Item(5, 12).Value = 1.033866968510776
$ws.Cells.Item(5, 13).Value = 1.055508005000181
$ws.Cells.Item(5, 14).Value = 1.037880980061265
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.03200114501155
$ws.Cells.Item(6, 4).Value = 1.042116927940622
$ws.Cells.Item(6, 5).Value = 1.031462212712136
$ws.Cells.Item(6, 6).Value = 1.053160710911756
$ws.Cells.Item(6, 9).Value = 1.037788558967006
$ws.Cells.Item(6, 10).Value = 1.036438276563449
$ws.Cells.Item(6, 11).Value = 1.04452543124654
$ws.Cells.Item(6, 12).Value = 1.033896931093152
$ws.Cells.Item(6, 13).Value = 1.05554264327686
$ws.Cells.Item(6, 14).Value = 1.037910137797301
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.031687719182785
$ws.Cells.Item(7, 4).Value = 1.041868290906854
$ws.Cells.Item(7, 5).Value = 1.031195417683745
$ws.Cells.Item(7, 6).Value = 1.052863046996678
$ws.Cells.Item(7, 9).Value = 1.037719062373543
$ws.Cells.Item(7, 10).Value = 1.036238039461134
$ws.Cells.Item(7, 11).Value = 1.044336818551302
$ws.Cells.Item(7, 12).Value = 1.033690899069269
$ws.Cells.Item(7, 13).Value = 1.055304422813343
$ws.Cells.Item(7, 14).Value = 1.037709616335333
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.030378157503618
$ws.Cells.Item(8, 4).Value = 1.040828713748269
$ws.Cells.Item(8, 5).Value = 1.030081473759513
$ws.Cells.Item(8, 6).Value = 1.051619108681457
$ws.Cells.Item(8, 9).Value = 1.037425099551362
$ws.Cells.Item(8, 10).Value = 1.035400405175381
$ws.Cells.Item(8, 11).Value = 1.043546805838569
$ws.Cells.Item(8, 12).Value = 1.032829633205019
$ws.Cells.Item(8, 13).Value = 1.054307682121452
$ws.Cells.Item(8, 14).Value = 1.036870792512817
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.028072486403439
$ws.Cells.Item(9, 4).Value = 1.0389957306675
$ws.Cells.Item(9, 5).Value = 1.028123211270656
$ws.Cells.Item(9, 6).Value = 1.049428136715575
$ws.Cells.Item(9, 9).Value = 1.03689402392522
$ws.Cells.Item(9, 10).Value = 1.033921884738602
$ws.Cells.Item(9, 11).Value = 1.042148566855685
$ws.Cells.Item(9, 12).Value = 1.03131170007743
$ws.Cells.Item(9, 13).Value = 1.052547523539786
$ws.Cells.Item(9, 14).Value = 1.035390172407429
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.026536705988622
$ws.Cells.Item(10, 4).Value = 1.037773055172883
$ws.Cells.Item(10, 5).Value = 1.026820886953248
$ws.Cells.Item(10, 6).Value = 1.047968226764966
$ws.Cells.Item(10, 9).Value = 1.036531214020423
$ws.Cells.Item(10, 10).Value = 1.032934562289053
$ws.Cells.Item(10, 11).Value = 1.041212330768326
$ws.Cells.Item(10, 12).Value = 1.030299605216342
$ws.Cells.Item(10, 13).Value = 1.051371602279449
$ws.Cells.Item(10, 14).Value = 1.034401447846754
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.025872017280631
$ws.Cells.Item(11, 4).Value = 1.037243473066521
$ws.Cells.Item(11, 5).Value = 1.026257733179794
$ws.Cells.Item(11, 6).Value = 1.0473362570096
$ws.Cells.Item(11, 9).Value = 1.036372042058834
$ws.Cells.Item(11, 10).Value = 1.032506659199191
$ws.Cells.Item(11, 11).Value = 1.040805972973202
$ws.Cells.Item(11, 12).Value = 1.029861332705592
$ws.Cells.Item(11, 13).Value = 1.050861837726298
$ws.Cells.Item(11, 14).Value = 1.033972937085423
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.025625169822423
$ws.Cells.Item(12, 4).Value = 1.03704674029037
$ws.Cells.Item(12, 5).Value = 1.026048668068286
$ws.Cells.Item(12, 6).Value = 1.047101543399691
$ws.Cells.Item(12, 9).Value = 1.036312607535184
$ws.Cells.Item(12, 10).Value = 1.03234765953171
$ws.Cells.Item(12, 11).Value = 1.040654889873658
$ws.Cells.Item(12, 12).Value = 1.029698535118273
$ws.Cells.Item(12, 13).Value = 1.05067240172211
$ws.Cells.Item(12, 14).Value = 1.033813711620177
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.025678117284858
$ws.Cells.Item(13, 4).Value = 1.037088941136947
$ws.Cells.Item(13, 5).Value = 1.026093508013479
$ws.Cells.Item(13, 6).Value = 1.047151888979661
$ws.Cells.Item(13, 9).Value = 1.036325370507005
$ws.Cells.Item(13, 10).Value = 1.032381768088678
$ws.Cells.Item(13, 11).Value = 1.040687304203598
$ws.Cells.Item(13, 12).Value = 1.02973345590859
$ws.Cells.Item(13, 13).Value = 1.050713040300639
$ws.Cells.Item(13, 14).Value = 1.033847868615208
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.025851611821552
$ws.Cells.Item(14, 4).Value = 1.037227211525294
$ws.Cells.Item(14, 5).Value = 1.02624044943567
$ws.Cells.Item(14, 6).Value = 1.047316854930937
$ws.Cells.Item(14, 9).Value = 1.03636713552363
$ws.Cells.Item(14, 10).Value = 1.032493517404837
$ws.Cells.Item(14, 11).Value = 1.040793487324802
$ws.Cells.Item(14, 12).Value = 1.02984787587716
$ws.Cells.Item(14, 13).Value = 1.050846180652077
$ws.Cells.Item(14, 14).Value = 1.033959776628214
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.025958513837826
$ws.Cells.Item(15, 4).Value = 1.037312401535394
$ws.Cells.Item(15, 5).Value = 1.026331000197566
$ws.Cells.Item(15, 6).Value = 1.04741849964847
$ws.Cells.Item(15, 9).Value = 1.036392827135867
$ws.Cells.Item(15, 10).Value = 1.032562362256911
$ws.Cells.Item(15, 11).Value = 1.040858891227165
$ws.Cells.Item(15, 12).Value = 1.029918373330423
$ws.Cells.Item(15, 13).Value = 1.050928201345792
$ws.Cells.Item(15, 14).Value = 1.034028719247874
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.026580825837243
$ws.Cells.Item(16, 4).Value = 1.037808198598585
$ws.Cells.Item(16, 5).Value = 1.026858277730478
$ws.Cells.Item(16, 6).Value = 1.048010172405945
$ws.Cells.Item(16, 9).Value = 1.036541734105156
$ws.Cells.Item(16, 10).Value = 1.032962952725614
$ws.Cells.Item(16, 11).Value = 1.041239279202746
$ws.Cells.Item(16, 12).Value = 1.030328691369772
$ws.Cells.Item(16, 13).Value = 1.05140542146332
$ws.Cells.Item(16, 14).Value = 1.034429878600991
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.026971270081556
$ws.Cells.Item(17, 4).Value = 1.038119158037199
$ws.Cells.Item(17, 5).Value = 1.027189229444569
$ws.Cells.Item(17, 6).Value = 1.048381362174042
$ws.Cells.Item(17, 9).Value = 1.036634584867682
$ws.Cells.Item(17, 10).Value = 1.03321412962092
$ws.Cells.Item(17, 11).Value = 1.041477629645745
$ws.Cells.Item(17, 12).Value = 1.030586065801348
$ws.Cells.Item(17, 13).Value = 1.051704613560073
$ws.Cells.Item(17, 14).Value = 1.034681412196299
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.027199039913416
$ws.Cells.Item(18, 4).Value = 1.038300520343947
$ws.Cells.Item(18, 5).Value = 1.027382341369187
$ws.Cells.Item(18, 6).Value = 1.048597888277313
$ws.Cells.Item(18, 9).Value = 1.036688543113353
$ws.Cells.Item(18, 10).Value = 1.033360599576804
$ws.Cells.Item(18, 11).Value = 1.041616562618429
$ws.Cells.Item(18, 12).Value = 1.03073618513423
$ws.Cells.Item(18, 13).Value = 1.051879070953206
$ws.Cells.Item(18, 14).Value = 1.034828090156321
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.027276708672073
$ws.Cells.Item(19, 4).Value = 1.038362357617618
$ws.Cells.Item(19, 5).Value = 1.027448199965497
$ws.Cells.Item(19, 6).Value = 1.048671721004869
$ws.Cells.Item(19, 9).Value = 1.03670690753769
$ws.Cells.Item(19, 10).Value = 1.033410535706109
$ws.Cells.Item(19, 11).Value = 1.041663919408027
$ws.Cells.Item(19, 12).Value = 1.030787371449184
$ws.Cells.Item(19, 13).Value = 1.051938546842257
$ws.Cells.Item(19, 14).Value = 1.034878097200657
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.026929375971931
$ws.Cells.Item(20, 4).Value = 1.038085796584317
$ws.Cells.Item(20, 5).Value = 1.027153713860173
$ws.Cells.Item(20, 6).Value = 1.048341535219213
$ws.Cells.Item(20, 9).Value = 1.036624643546624
$ws.Cells.Item(20, 10).Value = 1.033187184568726
$ws.Cells.Item(20, 11).Value = 1.041452066493274
$ws.Cells.Item(20, 12).Value = 1.030558452247343
$ws.Cells.Item(20, 13).Value = 1.051672518896467
$ws.Cells.Item(20, 14).Value = 1.03465442887904
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.025800520690476
$ws.Cells.Item(21, 4).Value = 1.037186494952892
$ws.Cells.Item(21, 5).Value = 1.026197175659861
$ws.Cells.Item(21, 6).Value = 1.047268275787012
$ws.Cells.Item(21, 9).Value = 1.03635484534784
$ws.Cells.Item(21, 10).Value = 1.032460611596607
$ws.Cells.Item(21, 11).Value = 1.040762223003828
$ws.Cells.Item(21, 12).Value = 1.029814182144098
$ws.Cells.Item(21, 13).Value = 1.050806976525223
$ws.Cells.Item(21, 14).Value = 1.033926824089961
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.025091039673188
$ws.Cells.Item(22, 4).Value = 1.036620939187171
$ws.Cells.Item(22, 5).Value = 1.025596429516652
$ws.Cells.Item(22, 6).Value = 1.046593637989802
$ws.Cells.Item(22, 9).Value = 1.0361834135598
$ws.Cells.Item(22, 10).Value = 1.032003454473855
$ws.Cells.Item(22, 11).Value = 1.040327659846969
$ws.Cells.Item(22, 12).Value = 1.02934620912286
$ws.Cells.Item(22, 13).Value = 1.050262274166825
$ws.Cells.Item(22, 14).Value = 1.033469017751658
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.025467122808535
$ws.Cells.Item(23, 4).Value = 1.036920762918247
$ws.Cells.Item(23, 5).Value = 1.025914832851427
$ws.Cells.Item(23, 6).Value = 1.0469512604937
$ws.Cells.Item(23, 9).Value = 1.036274463162343
$ws.Cells.Item(23, 10).Value = 1.032245833378226
$ws.Cells.Item(23, 11).Value = 1.040558108546093
$ws.Cells.Item(23, 12).Value = 1.029594292234682
$ws.Cells.Item(23, 13).Value = 1.05055107845821
$ws.Cells.Item(23, 14).Value = 1.033711740861875
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.026948306023024
$ws.Cells.Item(24, 4).Value = 1.03810087123546
$ws.Cells.Item(24, 5).Value = 1.027169761597296
$ws.Cells.Item(24, 6).Value = 1.04835953125205
$ws.Cells.Item(24, 9).Value = 1.036629136219763
$ws.Cells.Item(24, 10).Value = 1.033199359993152
$ws.Cells.Item(24, 11).Value = 1.041463617668294
$ws.Cells.Item(24, 12).Value = 1.030570929632509
$ws.Cells.Item(24, 13).Value = 1.051687021267405
$ws.Cells.Item(24, 14).Value = 1.034666621593966
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.028668324538738
$ws.Cells.Item(25, 4).Value = 1.039469726377461
$ws.Cells.Item(25, 5).Value = 1.028628912114196
$ws.Cells.Item(25, 6).Value = 1.049994429885613
$ws.Cells.Item(25, 9).Value = 1.037032865939394
$ws.Cells.Item(25, 10).Value = 1.034304410065386
$ws.Cells.Item(25, 11).Value = 1.042510767239823
$ws.Cells.Item(25, 12).Value = 1.031704149954891
$ws.Cells.Item(25, 13).Value = 1.053003008260767
$ws.Cells.Item(25, 14).Value = 1.035773240964053
